# The dataset workbook's first tab ("Sheet4") is renamed to "final" to
# mark it as the finished/consolidated mapping sheet, matching the
# "Reading files and mapping working." commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Name = "final"
